# DevHub Meeting Minutes - update sprint schedule content.
# The sheet's header columns are reordered (Meeting Type moves next to
# Sprint, ahead of Date/Time/Duration) and the six meeting rows are
# replaced with the new sprint planning / mid-sprint sync / review data.
# The underlying Table4 (A2:H8) stays the same shape; its column names
# simply follow whatever is written into row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (A1) ---
$ws.Range("A1").Value = "DevHub – Meeting Minutes Summary (Sprint 1 to 6)"

# --- Header row (row 2): Sprint, Meeting Type, Date, Time, Duration, Attendees, Agenda, Decisions ---
$ws.Range("A2").Value = "Sprint"
$ws.Range("B2").Value = "Meeting Type"
$ws.Range("C2").Value = "Date"
$ws.Range("D2").Value = "Time"
$ws.Range("E2").Value = "Duration"
$ws.Range("F2").Value = "Attendees"
$ws.Range("G2").Value = "Agenda / Discussion Points"
$ws.Range("H2").Value = "Decisions / Outcomes"
$ws.Rows.Item(2).RowHeight = 28.8

# --- Data rows 3-8 ---
# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Sprint Planning"
$ws.Range("C3").Value = "Wed, Mar 5, 2025"
$ws.Range("D3").Value = "04:00 PM – 04:45 PM"
$ws.Range("E3").Value = "45 mins"
$ws.Range("F3").Value = "Fatima, Rushba, Saad"
$ws.Range("G3").Value = "GitHub setup, schema delays, branching confusion"
$ws.Range("H3").Value = "Branching policy set, schema finalized"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Sprint Planning"
$ws.Range("C4").Value = "Wed, Mar 12, 2025"
$ws.Range("D4").Value = "03:30 PM – 04:15 PM"
$ws.Range("E4").Value = "45 mins"
$ws.Range("F4").Value = "Fatima, Rushba, Saad"
$ws.Range("G4").Value = "UI planning, SEO basics, Figma templates usage"
$ws.Range("H4").Value = "Tasks divided, SEO learning assigned"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Mid-Sprint Sync"
$ws.Range("C5").Value = "Wed, Mar 19, 2025"
$ws.Range("D5").Value = "04:00 PM – 04:30 PM"
$ws.Range("E5").Value = "30 mins"
$ws.Range("F5").Value = "Full Team"
$ws.Range("G5").Value = "Backend CRUD status, integration, MongoDB design"
$ws.Range("H5").Value = "Prioritized frontend/backend sync, deferred docs"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Mid-Sprint Sync"
$ws.Range("C6").Value = "Wed, Mar 26, 2025"
$ws.Range("D6").Value = "03:00 PM – 03:30 PM"
$ws.Range("E6").Value = "30 mins"
$ws.Range("F6").Value = "Full Team"
$ws.Range("G6").Value = "Rating system, email reliability, social API limits"
$ws.Range("H6").Value = "Added spam control, email testing, limited APIs"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Sprint Planning"
$ws.Range("C7").Value = "Wed, Apr 2, 2025"
$ws.Range("D7").Value = "04:15 PM – 05:00 PM"
$ws.Range("E7").Value = "45 mins"
$ws.Range("F7").Value = "Full Team"
$ws.Range("G7").Value = "Plan for quizzes, progress dashboard, interactive snippets"
$ws.Range("H7").Value = "Assigned features module-wise, used schema templates"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Sprint Review Prep"
$ws.Range("C8").Value = "Wed, Apr 9, 2025"
$ws.Range("D8").Value = "05:00 PM – 05:30 PM"
$ws.Range("E8").Value = "30 mins"
$ws.Range("F8").Value = "Full Team"
$ws.Range("G8").Value = "Analytics integration, accessibility rechecks, performance fix"
$ws.Range("H8").Value = "Confirmed tracking setup, screenshots finalized"

# --- Selection ends on B2, matching the saved workbook view ---
$ws.Range("B2").Select()
